$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.Value = "'" + $val
    $c.Style = $origStyle
}

Set-TextValue 'D2' '69.075.04'
$ws.Range('E2').Value = '  -0.59%  '

Set-TextValue 'D3' '2.469.15'
$ws.Range('E3').Value = '  -0.94%  '

Set-TextValue 'D5' '558.72'
$ws.Range('E5').Value = '  -1.95%  '

Set-TextValue 'D6' '162.89'
$ws.Range('E6').Value = '  -2.20%  '

Set-TextValue 'D8' '0.505'
$ws.Range('E8').Value = '  -0.96%  '

Set-TextValue 'D9' '2.468.20'
$ws.Range('E9').Value = '  -0.93%  '

Set-TextValue 'D10' '0.151'
$ws.Range('E10').Value = '  -5.12%  '

$ws.Range('E11').Value = '  -0.55%  '

Set-TextValue 'D12' '0.335'
$ws.Range('E12').Value = '  -4.65%  '

Set-TextValue 'D13' '4.83'
$ws.Range('E13').Value = '  -0.68%  '

Set-TextValue 'D14' '2.920.76'
$ws.Range('E14').Value = '  -1.03%  '

Set-TextValue 'D15' '68.978.13'
$ws.Range('E15').Value = '  -0.68%  '

Set-TextValue 'D16' '0.0000170'
$ws.Range('E16').Value = '  -3.02%  '

Set-TextValue 'D17' '23.62'
$ws.Range('E17').Value = '  -2.49%  '

Set-TextValue 'D18' '2.455.04'
$ws.Range('E18').Value = '  -2.28%  '

Set-TextValue 'D19' '10.77'
$ws.Range('E19').Value = '  -4.14%  '

Set-TextValue 'D20' '342.66'
$ws.Range('E20').Value = '  -3.10%  '

Set-TextValue 'D21' '7.05'
$ws.Range('E21').Value = '  -4.39%  '

Set-TextValue 'D22' '3.81'
$ws.Range('E22').Value = '  -2.48%  '

$ws.Range('E23').Value = '  -0.47%  '

Set-TextValue 'D24' '1.01'
$ws.Range('E24').Value = '  +0.57%  '

Set-TextValue 'D25' '1.89'
$ws.Range('E25').Value = '  -1.22%  '

Set-TextValue 'D26' '67.02'
$ws.Range('E26').Value = '  -3.37%  '

Set-TextValue 'D27' '3.69'
$ws.Range('E27').Value = '  -2.75%  '

Set-TextValue 'D28' '2.597.36'
$ws.Range('E28').Value = '  -1.01%  '

Set-TextValue 'D29' '0.999'
$ws.Range('E29').Value = '  -0.06%  '

Set-TextValue 'D30' '8.16'
$ws.Range('E30').Value = '  -5.56%  '

Set-TextValue 'D31' '0.0₃0821'
$ws.Range('E31').Value = '  -6.05%  '

Set-TextValue 'D32' '7.19'
$ws.Range('E32').Value = '  -5.28%  '

Set-TextValue 'D33' '438.37'
$ws.Range('E33').Value = '  -1.07%  '

Set-TextValue 'D34' '0.999'
$ws.Range('E34').Value = '  -0.04%  '

Set-TextValue 'D35' '1.15'
$ws.Range('E35').Value = '  -3.98%  '

$ws.Range('E36').Value = '  -5.67%  '

Set-TextValue 'D37' '157.32'
$ws.Range('E37').Value = '  +1.99%  '

$ws.Range('E38').Value = '  -0.02%  '

$ws.Range('E39').Value = '  +0.04%  '

$ws.Range('E40').Value = '  -3.33%  '

Set-TextValue 'D41' '17.90'
$ws.Range('E41').Value = '  -1.31%  '

Set-TextValue 'D42' '0.304'
$ws.Range('E42').Value = '  -3.24%  '

Set-TextValue 'D43' '4.46'
$ws.Range('E43').Value = '  -2.90%  '

Set-TextValue 'D44' '37.47'
$ws.Range('E44').Value = '  -0.79%  '

Set-TextValue 'D45' '1.48'
$ws.Range('E45').Value = '  -6.24%  '

$ws.Range('E46').Value = '  +2.94%  '

Set-TextValue 'D47' '2.08'
$ws.Range('E47').Value = '  -4.99%  '

Set-TextValue 'D48' '133.43'
$ws.Range('E48').Value = '  -3.77%  '

Set-TextValue 'D49' '3.36'
$ws.Range('E49').Value = '  -2.18%  '

Set-TextValue 'D50' '0.0718'
$ws.Range('E50').Value = '  -0.56%  '

Set-TextValue 'D51' '0.485'
$ws.Range('E51').Value = '  -4.05%  '
